$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.436.65"
$ws.Range("E2").Value = "  -2.22%  "
$ws.Range("D3").Value = "2.522.00"
$ws.Range("E3").Value = "  -5.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.96"
$ws.Range("E5").Value = "  -3.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.31"
$ws.Range("E6").Value = "  -3.06%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -2.84%  "
$ws.Range("D9").Value = "2.524.22"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  -4.60%  "
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("E12").Value = "  -3.77%  "
$ws.Range("E13").Value = "  -4.13%  "
$ws.Range("D14").Value = "2.990.18"
$ws.Range("E14").Value = "  -5.05%  "
$ws.Range("D15").Value = "70.303.22"
$ws.Range("E15").Value = "  -2.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000179"
$ws.Range("E16").Value = "  -3.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.82"
$ws.Range("E17").Value = "  -5.47%  "
$ws.Range("D18").Value = "2.525.38"
$ws.Range("E18").Value = "  -5.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.46"
$ws.Range("E19").Value = "  -6.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.53"
$ws.Range("E20").Value = "  -8.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "354.15"
$ws.Range("E21").Value = "  -4.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.93"
$ws.Range("E22").Value = "  -5.65%  "
$ws.Range("E23").Value = "  -3.80%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.20"
$ws.Range("E25").Value = "  -4.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.04"
$ws.Range("E26").Value = "  -6.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.18"
$ws.Range("E27").Value = "  -5.79%  "
$ws.Range("D28").Value = "2.656.96"
$ws.Range("E28").Value = "  -5.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.04"
$ws.Range("E29").Value = "  +4.14%  "
$ws.Range("D30").Value = "0.0₃0907"
$ws.Range("E30").Value = "  -6.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.83"
$ws.Range("E31").Value = "  -2.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "477.72"
$ws.Range("E32").Value = "  -4.59%  "
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.75"
$ws.Range("E34").Value = "  -3.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.27"
$ws.Range("E36").Value = "  -3.55%  "
$ws.Range("E37").Value = "  +3.38%  "
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.84"
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.57"
$ws.Range("E39").Value = "  -4.98%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.30"
$ws.Range("E41").Value = "  -5.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.64"
$ws.Range("E42").Value = "  -6.94%  "
$ws.Range("E43").Value = "  -3.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.69"
$ws.Range("E44").Value = "  -5.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("E45").Value = "  -6.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.29"
$ws.Range("E46").Value = "  -3.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.43"
$ws.Range("E47").Value = "  -9.57%  "
$ws.Range("E48").Value = "  -5.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.522"
$ws.Range("E49").Value = "  -6.84%  "
$ws.Range("E50").Value = "  -7.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.595"
$ws.Range("E51").Value = "  -1.78%  "
